$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the date values in column F (rows 2-7), shifting each date
# forward by 552 days while preserving existing cell formatting.
$ws.Range("F2").Value = 45957
$ws.Range("F3").Value = 45956
$ws.Range("F4").Value = 45955
$ws.Range("F5").Value = 45954
$ws.Range("F6").Value = 45953
$ws.Range("F7").Value = 45952
